# Fruta / hortaliza, semanal
# Insert a new weekly record at row 143, pushing the existing
# rows 143-169 down to 144-170.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 143 (shifts 143:169 -> 144:170)
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new weekly observation
$ws.Cells.Item(143, 1).Value = 11
$ws.Cells.Item(143, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(143, 3).Value = "Bíobío"
$ws.Cells.Item(143, 4).Value = 45015
$ws.Cells.Item(143, 5).Value = 8
$ws.Cells.Item(143, 6).Value = 100112024
$ws.Cells.Item(143, 7).Value = "Choclo"
$ws.Cells.Item(143, 8).Value = "Choclero"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 10000
$ws.Cells.Item(143, 11).Value = 300
$ws.Cells.Item(143, 12).Value = 400
$ws.Cells.Item(143, 13).Value = 350
$ws.Cells.Item(143, 14).Value = "$/unidad"
$ws.Cells.Item(143, 15).Value = "Región Metropolitana"
$ws.Cells.Item(143, 16).Value = 350
$ws.Cells.Item(143, 17).Value = 1
$ws.Cells.Item(143, 18).Value = "Hortaliza"
